$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (21:04 -> 22:04)
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 22:04"

# Update country rows affected by re-sorting (Peru/India, Guinea/Bolivia,
# Mali/Maldivas, Republica del Chad insertion that shifts several rows down,
# Belice/Nueva Caledonia, Montserrat/Seychelles) together with the refreshed
# daily case statistics for each impacted row.

$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1283951
$ws.Range("C4").Value = 20859
$ws.Range("D4").Value = 214844
$ws.Range("E4").Value = 992657
$ws.Range("F4").Value = 16943
$ws.Range("G4").Value = 1651
$ws.Range("H4").Value = 76450

$ws.Range("A10").Value = "Alemania"
$ws.Range("B10").Value = 169175
$ws.Range("C10").Value = 1013
$ws.Range("D10").Value = 139900
$ws.Range("E10").Value = 21918
$ws.Range("F10").Value = 1823
$ws.Range("G10").Value = 82
$ws.Range("H10").Value = 7357

$ws.Range("A16").Value = "Peru"
$ws.Range("B16").Value = 58526
$ws.Range("C16").Value = 3709
$ws.Range("D16").Value = 18388
$ws.Range("E16").Value = 38511
$ws.Range("F16").Value = 722
$ws.Range("G16").Value = 94
$ws.Range("H16").Value = 1627

$ws.Range("A17").Value = "India"
$ws.Range("B17").Value = 56351
$ws.Range("C17").Value = 3364
$ws.Range("D17").Value = 16776
$ws.Range("E17").Value = 37686
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 104
$ws.Range("H17").Value = 1889

$ws.Range("A22").Value = "Suiza"
$ws.Range("B22").Value = 30126
$ws.Range("C22").Value = 66
$ws.Range("D22").Value = 25900
$ws.Range("E22").Value = 2416
$ws.Range("F22").Value = 121
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 1810

$ws.Range("A37").Value = "Rumania"
$ws.Range("B37").Value = 14499
$ws.Range("C37").Value = 392
$ws.Range("D37").Value = 6144
$ws.Range("E37").Value = 7467
$ws.Range("F37").Value = 234
$ws.Range("G37").Value = 24
$ws.Range("H37").Value = 888

$ws.Range("A72").Value = "Uzbekistan"
$ws.Range("B72").Value = 2298
$ws.Range("C72").Value = 65
$ws.Range("D72").Value = 1656
$ws.Range("E72").Value = 632
$ws.Range("F72").Value = 8
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 10

$ws.Range("A77").Value = "Guinea"
$ws.Range("B77").Value = 1927
$ws.Range("C77").Value = 71
$ws.Range("D77").Value = 629
$ws.Range("E77").Value = 1287
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 11

$ws.Range("A78").Value = "Bolivia"
$ws.Range("B78").Value = 1886
$ws.Range("C78").Value = 84
$ws.Range("D78").Value = 198
$ws.Range("E78").Value = 1597
$ws.Range("F78").Value = 3
$ws.Range("G78").Value = 5
$ws.Range("H78").Value = 91

$ws.Range("A102").Value = "Sri Lanka"
$ws.Range("B102").Value = 823
$ws.Range("C102").Value = 26
$ws.Range("D102").Value = 232
$ws.Range("E102").Value = 582
$ws.Range("F102").Value = 1
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 9

$ws.Range("A106").Value = "Costa Rica"
$ws.Range("B106").Value = 765
$ws.Range("C106").Value = 4
$ws.Range("D106").Value = 445
$ws.Range("E106").Value = 314
$ws.Range("F106").Value = 6
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 6

$ws.Range("A112").Value = "Mali"
$ws.Range("B112").Value = 650
$ws.Range("C112").Value = 19
$ws.Range("D112").Value = 271
$ws.Range("E112").Value = 347
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 32

$ws.Range("A113").Value = "Maldivas"
$ws.Range("B113").Value = 648
$ws.Range("C113").Value = 31
$ws.Range("D113").Value = 20
$ws.Range("E113").Value = 625
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 3

$ws.Range("A129").Value = "Estado de Palestina"
$ws.Range("B129").Value = 375
$ws.Range("C129").Value = 1
$ws.Range("D129").Value = 176
$ws.Range("E129").Value = 197
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 2

$ws.Range("A136").Value = "Republica del Chad"
$ws.Range("B136").Value = 253
$ws.Range("C136").Value = 83
$ws.Range("D136").Value = 50
$ws.Range("E136").Value = 176
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 10
$ws.Range("H136").Value = 27

$ws.Range("A137").Value = "Sierra Leona"
$ws.Range("B137").Value = 231
$ws.Range("C137").Value = 6
$ws.Range("D137").Value = 54
$ws.Range("E137").Value = 161
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = 16

$ws.Range("A138").Value = "Cabo Verde"
$ws.Range("B138").Value = 218
$ws.Range("C138").Value = 27
$ws.Range("D138").Value = 38
$ws.Range("E138").Value = 178
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 2

$ws.Range("A139").Value = "Madagascar"
$ws.Range("B139").Value = 193
$ws.Range("C139").Value = 35
$ws.Range("D139").Value = 101
$ws.Range("E139").Value = 92
$ws.Range("F139").Value = 1
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0

$ws.Range("A140").Value = "Etiopia"
$ws.Range("B140").Value = 191
$ws.Range("C140").Value = 29
$ws.Range("D140").Value = 93
$ws.Range("E140").Value = 94
$ws.Range("F140").Value = 1
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 4

$ws.Range("A141").Value = "Islas Feroe"
$ws.Range("B141").Value = 187
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 185
$ws.Range("E141").Value = 2
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 0

$ws.Range("A142").Value = "Martinica"
$ws.Range("B142").Value = 182
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 83
$ws.Range("E142").Value = 85
$ws.Range("F142").Value = 5
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 14

$ws.Range("A143").Value = "Liberia"
$ws.Range("B143").Value = 178
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 75
$ws.Range("E143").Value = 83
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 20

$ws.Range("A144").Value = "Birmania"
$ws.Range("B144").Value = 176
$ws.Range("C144").Value = 15
$ws.Range("D144").Value = 62
$ws.Range("E144").Value = 108
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 6

$ws.Range("A145").Value = "Santo Tome y Principe"
$ws.Range("B145").Value = 174
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 4
$ws.Range("E145").Value = 167
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 3

$ws.Range("A153").Value = "Togo"
$ws.Range("B153").Value = 135
$ws.Range("C153").Value = 7
$ws.Range("D153").Value = 85
$ws.Range("E153").Value = 41
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 9

$ws.Range("A191").Value = "Belice"
$ws.Range("B191").Value = 18
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 16
$ws.Range("E191").Value = 0
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 2

$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("B192").Value = 18
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 18
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

$ws.Range("A205").Value = "Montserrat"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 7
$ws.Range("E205").Value = 3
$ws.Range("F205").Value = 1
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1

$ws.Range("A206").Value = "Seychelles"
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 8
$ws.Range("E206").Value = 3
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
